$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOB1036: Geometria Analitica
# (Requisito fraco)" requirement line - this paragraph (and everything
# before it) must stay untouched.
$rng = $d.Content
$rng.Find.Execute("LOB1036: Geometria Anal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hitStart = $rng.Start

$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $hitStart -and $hitStart -lt $p.Range.End) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    # Walk forward from the paragraph right after the anchor, looking for
    # the stray empty paragraph that forces a page break (the leftover
    # "Ver no Jupiter / Salvar em pdf / Salvar em docx" block ends with
    # it). Everything from the paragraph right after the anchor through
    # that page-break paragraph (inclusive) is removed.
    $scanIndex = $anchorIndex + 1
    $deleteEndIndex = -1
    while ($scanIndex -le $d.Paragraphs.Count) {
        $cand = $d.Paragraphs.Item($scanIndex)
        if ($cand.Format.PageBreakBefore) {
            $deleteEndIndex = $scanIndex
            break
        }
        $scanIndex = $scanIndex + 1
    }

    if ($deleteEndIndex -gt 0) {
        $anchorPara = $d.Paragraphs.Item($anchorIndex)
        $deleteEndPara = $d.Paragraphs.Item($deleteEndIndex)
        $delRange = $d.Range($anchorPara.Range.End, $deleteEndPara.Range.End)
        $delRange.Delete()
    }
}
